$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3 and each language sheet's Status column C2:C3)
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2) Overview columns E/F widen (text got longer)
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# 3) zh-cn sheet: fill in "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime", add hyperlink on I (Latest Target File),
#    widen columns C (Status) and J (Latest Handback File)
# ---------------------------------------------------------------------------
$zh.Range("I2").Value = "a.md"
$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-01 16:44:39"

$zh.Range("I3").Value = "a.md"
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-01 16:44:39"

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/b.md", $null, $null, "b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")

# ---------------------------------------------------------------------------
# 4) de-de sheet: same shape of change, but with de-de filenames/datetime
# ---------------------------------------------------------------------------
$de.Range("I2").Value = "a.md"
$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K2").Value = "2016-09-01 16:44:46"

$de.Range("I3").Value = "a.md"
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K3").Value = "2016-09-01 16:44:46"

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/b.md", $null, $null, "b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md", $null, $null, "a.md")
